$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Marking -> Right (B11) 6 -> 9, Wrong (C11) 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12: Total -> Right (B12) 114 -> 171, Max text (E12) "114/168" -> "171/252"
$ws.Range("B12").Value = 171
$ws.Range("E12").Value = "171/252"
